$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LTER_sites")
$wsDatasets = $wb.Worksheets.Item("datasets")

# New site record: Collelongo-Selva Piana (row 4)
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "9b1d144a-dc37-4b0e-8cda-1dda1d7667da"
$ws.Range("C4").Value = "Collelongo-Selva Piana "
$ws.Range("D4").Value = "Terrestrial"
$ws.Range("D4").Font.Bold = $true
$ws.Range("E4").Formula = "=TRUE"
$ws.Range("F4").Value = "https://deims.org/9b1d144a-dc37-4b0e-8cda-1dda1d7667da"
$ws.Range("G4").Value = "collelongo"

# Match the bold + left aligned style already used by the "active" column
$ws.Range("E2").Copy()
$ws.Range("E4").PasteSpecial(-4122)

# Update the active sheet/selection state to reflect the latest work:
# "datasets" was being reviewed around D26, then focus moved back to
# "LTER_sites" to review the newly added row near C15.
$wsDatasets.Activate()
$wsDatasets.Range("D26").Select()

$ws.Activate()
$ws.Range("C15").Select()
